$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 88.071501905872964
$ws.Range("D2").Value = 109.386817814227
$ws.Range("E2").Value = 3009
$ws.Range("G2").Value = 377.384521484375
$ws.Range("H2").Value = 7110.14306640625

$ws.Range("C3").Value = 94.976976093559728
$ws.Range("D3").Value = 117.96357444086199
$ws.Range("E3").Value = 3009
$ws.Range("G3").Value = 451.80047607421875
$ws.Range("H3").Value = 8847.2685546875

$ws.Range("C4").Value = 14.018136161076503
$ws.Range("D4").Value = 17.410845444581959
$ws.Range("E4").Value = 3009
$ws.Range("G4").Value = 61.808502197265625
$ws.Range("H4").Value = 609.37957763671875

$ws.Range("C5").Value = 1.5112757105968371
$ws.Range("D5").Value = 1.8770389899180633
$ws.Range("E5").Value = 3009
$ws.Range("G5").Value = 5.3433041572570801
$ws.Range("H5").Value = 270.29360961914063

$ws.Range("C6").Value = 45.010208802747741
$ws.Range("D6").Value = 55.903707537868165
$ws.Range("E6").Value = 3009
$ws.Range("G6").Value = 76.867599487304688
$ws.Range("H6").Value = 5031.33349609375

$ws.Range("C7").Value = 28.417130347936801
$ws.Range("D7").Value = 35.294724975007838
$ws.Range("E7").Value = 3009
$ws.Range("G7").Value = 69.354133605957031
$ws.Range("H7").Value = 7411.89208984375

$ws.Range("C8").Value = 9.8233872607285964
$ws.Range("D8").Value = 12.200871343666551
$ws.Range("E8").Value = 3009
$ws.Range("G8").Value = 22.358097076416016
$ws.Range("H8").Value = 1172.123046875

$ws.Range("C9").Value = 18.268528258614303
$ws.Range("D9").Value = 22.689929412583172
$ws.Range("E9").Value = 3009
$ws.Range("G9").Value = 39.480476379394531
$ws.Range("H9").Value = 3918.30224609375

$ws.Range("C10").Value = 26.817168198696756
$ws.Range("D10").Value = 33.307535595575338
$ws.Range("E10").Value = 3009
$ws.Range("G10").Value = 28.311405181884766
$ws.Range("H10").Value = 2664.602783203125

$ws.Range("C11").Value = 5.3506149702113346
$ws.Range("D11").Value = 6.6455859827337651
$ws.Range("E11").Value = 3009
$ws.Range("G11").Value = 28.476335525512695
$ws.Range("H11").Value = 1993.67578125

$ws.Range("C12").Value = 20.209372464348288
$ws.Range("D12").Value = 25.100502154881791
$ws.Range("E12").Value = 3009
$ws.Range("G12").Value = 92.495353698730469
$ws.Range("H12").Value = 6024.12060546875

$ws.Range("C13").Value = 16.742273180458625
$ws.Range("D13").Value = 20.794285805291139
$ws.Range("E13").Value = 3009
$ws.Range("G13").Value = 30.151714324951172
$ws.Range("H13").Value = 3432.70751953125

$ws.Range("C14").Value = 19.665394595351557
$ws.Range("D14").Value = 24.424869165146379
$ws.Range("E14").Value = 3009
$ws.Range("G14").Value = 219.82382202148438
$ws.Range("H14").Value = 3318.596435546875

$ws.Range("C15").Value = 5.620282138927152
$ws.Range("D15").Value = 6.9805187861283144
$ws.Range("E15").Value = 3009
$ws.Range("G15").Value = 40.766231536865234
$ws.Range("H15").Value = 1814.9349365234375

$ws.Range("C16").Value = 28.7704038673875
$ws.Range("D16").Value = 35.733498839216402
$ws.Range("E16").Value = 3009
$ws.Range("G16").Value = 40.512855529785156
$ws.Range("H16").Value = 893.33746337890625

$ws.Range("C17").Value = 16.781085621918425
$ws.Range("D17").Value = 20.842491949860459
$ws.Range("E17").Value = 3009
$ws.Range("G17").Value = 34.360336303710938
$ws.Range("H17").Value = 729.48724365234375

$ws.Range("C18").Value = 8.7791864503535297
$ws.Range("D18").Value = 10.903950118828391
$ws.Range("E18").Value = 3009
$ws.Range("G18").Value = 2.5079085826873779
$ws.Range("H18").Value = 872.31597900390625

$ws.Range("C19").Value = 9.9464465566313756
$ws.Range("D19").Value = 12.353713879834064
$ws.Range("E19").Value = 3009
$ws.Range("G19").Value = 4.3238000869750977
$ws.Range("H19").Value = 741.22283935546875

$ws.Range("C20").Value = 13.179509357821312
$ws.Range("D20").Value = 16.369251651199761
$ws.Range("E20").Value = 3009
$ws.Range("G20").Value = 6.5477008819580078
$ws.Range("H20").Value = 982.15509033203125

$ws.Range("C21").Value = 6.6915920387672481
$ws.Range("D21").Value = 8.3111101554928926
$ws.Range("E21").Value = 3009
$ws.Range("G21").Value = 4.695777416229248
$ws.Range("H21").Value = 831.11102294921875

$ws.Range("C22").Value = 16.975265168240394
$ws.Range("D22").Value = 21.083666997004997
$ws.Range("E22").Value = 3009
$ws.Range("G22").Value = 28.486377716064453
$ws.Range("H22").Value = 1992.406494140625

$ws.Range("C23").Value = 7.4372123992098409
$ws.Range("D23").Value = 9.2371876337235843
$ws.Range("E23").Value = 3009
$ws.Range("G23").Value = 37.48876953125
$ws.Range("H23").Value = 720.5006103515625

$ws.Range("C24").Value = 17.566349087224207
$ws.Range("D24").Value = 21.817806737056085
$ws.Range("E24").Value = 3009
$ws.Range("G24").Value = 55.526317596435547
$ws.Range("H24").Value = 1745.424560546875

$ws.Range("C25").Value = 33.815221819460845
$ws.Range("D25").Value = 41.999278094214674
$ws.Range("E25").Value = 3009
$ws.Range("G25").Value = 76.578681945800781
$ws.Range("H25").Value = 839.98553466796875

$ws.Range("C26").Value = 13.281584191060851
$ws.Range("D26").Value = 16.496030918424331
$ws.Range("E26").Value = 3009
$ws.Range("G26").Value = 20.2012939453125
$ws.Range("H26").Value = 1154.72216796875

$ws.Range("C27").Value = 29.394199233785663
$ws.Range("D27").Value = 36.50826673041788
$ws.Range("E27").Value = 3009
$ws.Range("G27").Value = 143.6600341796875
$ws.Range("H27").Value = 2738.1201171875

$ws.Range("C28").Value = 11.080093315171419
$ws.Range("D28").Value = 13.761728929815677
$ws.Range("E28").Value = 3009
$ws.Range("G28").Value = 0.68808645009994507
$ws.Range("H28").Value = 3302.81494140625

$ws.Range("C29").Value = 17.59972508986078
$ws.Range("D29").Value = 21.859260514234933
$ws.Range("E29").Value = 3009
$ws.Range("H29").Value = 819.7222900390625
